# Auto-generated edit script applying scheduled data-refresh updates
# to the per-job H..N (price/profit) columns across all 8 class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H96").Value = 851.5833
$ws.Range("I96").Value = 569.5
$ws.Range("J96").Value = 1133.6666
$ws.Range("K96").Value = 1708.5
$ws.Range("L96").Value = 3400.9998
$ws.Range("M96").Value = -335.5
$ws.Range("N96").Value = -6146.9998

$ws.Range("H132").Value = 6212
$ws.Range("I132").Value = 6346.8887
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 19040.6661
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -16510.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1974.7037
$ws.Range("I32").Value = 1932.68
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 1932.68
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -1645.68
$ws.Range("N32").Value = -3074

$ws.Range("H45").Value = 2999.4
$ws.Range("I45").Value = 2650
$ws.Range("J45").Value = 3232.3333
$ws.Range("K45").Value = 2650
$ws.Range("L45").Value = 3232.3333
$ws.Range("M45").Value = -2273

$ws.Range("H80").Value = 85994.336
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 85994.336
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 85994.336
$ws.Range("N80").Value = -87990.336

$ws.Range("H83").Value = 85994.336
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 85994.336
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 257983.008
$ws.Range("N83").Value = -267967.008

$ws.Range("H122").Value = 1549.5714
$ws.Range("I122").Value = 1549.5714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4648.7142
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2198.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 408
$ws.Range("I94").Value = 408
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 408
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 43

$ws.Range("H105").Value = 3742.5334
$ws.Range("I105").Value = 2449.8333
$ws.Range("J105").Value = 4604.3335
$ws.Range("K105").Value = 2449.8333
$ws.Range("L105").Value = 4604.3335
$ws.Range("M105").Value = -702.8332999999998
$ws.Range("N105").Value = -8098.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5352.5356
$ws.Range("I22").Value = 1323.9286
$ws.Range("J22").Value = 9381.143
$ws.Range("K22").Value = 1323.9286
$ws.Range("L22").Value = 9381.143
$ws.Range("M22").Value = -973.9286
$ws.Range("N22").Value = -10081.143

$ws.Range("H41").Value = 13037.25
$ws.Range("I41").Value = 9049.666999999999
$ws.Range("J41").Value = 25000
$ws.Range("K41").Value = 9049.666999999999
$ws.Range("L41").Value = 25000
$ws.Range("M41").Value = -8621.666999999999
$ws.Range("N41").Value = -25856

$ws.Range("H86").Value = 4999.3335
$ws.Range("I86").Value = 4999.3335
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4999.3335
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3876.3335

$ws.Range("H89").Value = 4999.3335
$ws.Range("I89").Value = 4999.3335
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 24996.6675
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -19380.6675

$ws.Range("H132").Value = 1540.8
$ws.Range("I132").Value = 1355.5
$ws.Range("J132").Value = 1818.75
$ws.Range("K132").Value = 4066.5
$ws.Range("L132").Value = 5456.25
$ws.Range("M132").Value = -1536.5
$ws.Range("N132").Value = -10516.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 48
$ws.Range("I12").Value = 30
$ws.Range("J12").Value = 62.4
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 187.2
$ws.Range("M12").Value = 83
$ws.Range("N12").Value = -533.2

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H121").Value = 756.46155
$ws.Range("I121").Value = 818.2
$ws.Range("J121").Value = 717.875
$ws.Range("K121").Value = 2454.6
$ws.Range("L121").Value = 2153.625
$ws.Range("M121").Value = -1144.6
$ws.Range("N121").Value = -4773.625

$ws.Range("H139").Value = 3299.889
$ws.Range("I139").Value = 3539.8
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 10619.4
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -5479.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 345.57144
$ws.Range("I2").Value = 364
$ws.Range("J2").Value = 299.5
$ws.Range("K2").Value = 364
$ws.Range("L2").Value = 299.5
$ws.Range("M2").Value = -251
$ws.Range("N2").Value = -525.5

$ws.Range("H11").Value = 28708570
$ws.Range("I11").Value = 100250000
$ws.Range("J11").Value = 91999
$ws.Range("K11").Value = 100250000
$ws.Range("L11").Value = 91999
$ws.Range("M11").Value = -100249861

$ws.Range("H68").Value = 55000
$ws.Range("I68").Value = 55000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 55000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -54189

$ws.Range("H70").Value = 6409.1055
$ws.Range("I70").Value = 4658.6665
$ws.Range("J70").Value = 7984.5
$ws.Range("K70").Value = 4658.6665
$ws.Range("L70").Value = 7984.5
$ws.Range("M70").Value = -4388.6665

$ws.Range("H71").Value = 55000
$ws.Range("I71").Value = 55000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 165000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -160944

$ws.Range("H73").Value = 6409.1055
$ws.Range("I73").Value = 4658.6665
$ws.Range("J73").Value = 7984.5
$ws.Range("K73").Value = 4658.6665
$ws.Range("L73").Value = 7984.5
$ws.Range("M73").Value = -3722.6665

$ws.Range("H113").Value = 2961.125
$ws.Range("I113").Value = 2481.9167
$ws.Range("J113").Value = 4398.75
$ws.Range("K113").Value = 2481.9167
$ws.Range("L113").Value = 4398.75
$ws.Range("M113").Value = -311.9167000000002
$ws.Range("N113").Value = -8738.75

$ws.Range("H122").Value = 2369.5
$ws.Range("I122").Value = 1625.2142
$ws.Range("J122").Value = 4974.5
$ws.Range("K122").Value = 4875.642599999999
$ws.Range("L122").Value = 14923.5
$ws.Range("M122").Value = -2425.642599999999
$ws.Range("N122").Value = -19823.5

$ws.Range("H132").Value = 4227.3477
$ws.Range("I132").Value = 4100.125
$ws.Range("J132").Value = 4518.143
$ws.Range("K132").Value = 12300.375
$ws.Range("L132").Value = 13554.429
$ws.Range("M132").Value = -9770.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 505.07693
$ws.Range("I16").Value = 558.6
$ws.Range("J16").Value = 326.66666
$ws.Range("K16").Value = 558.6
$ws.Range("L16").Value = 326.66666
$ws.Range("M16").Value = -388.6
$ws.Range("N16").Value = -666.66666

$ws.Range("H22").Value = 1730.0454
$ws.Range("I22").Value = 1466.4445
$ws.Range("J22").Value = 1912.5385
$ws.Range("K22").Value = 1466.4445
$ws.Range("L22").Value = 1912.5385
$ws.Range("M22").Value = -1171.4445
$ws.Range("N22").Value = -2502.5385

$ws.Range("H27").Value = 1730.0454
$ws.Range("I27").Value = 1466.4445
$ws.Range("J27").Value = 1912.5385
$ws.Range("K27").Value = 1466.4445
$ws.Range("L27").Value = 1912.5385
$ws.Range("M27").Value = -1359.4445
$ws.Range("N27").Value = -2126.5385

$ws.Range("H46").Value = 2349.9375
$ws.Range("I46").Value = 2222.8462
$ws.Range("J46").Value = 2900.6667
$ws.Range("K46").Value = 2222.8462
$ws.Range("L46").Value = 2900.6667
$ws.Range("M46").Value = -2034.8462

$ws.Range("H98").Value = 58947
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 58947
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 58947
$ws.Range("N98").Value = -64937

$ws.Range("H122").Value = 5306.222
$ws.Range("I122").Value = 5306.222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15918.666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13468.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 4103.8
$ws.Range("I23").Value = 4504.75
$ws.Range("J23").Value = 2500
$ws.Range("K23").Value = 4504.75
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = -4275.75
$ws.Range("N23").Value = -2958
